# Update the "pie chart" column header in the report to a plain text label.
# Previously C5 held a rich-text shared string "Thời gian tham gia" + bold " (ngày)".
# Now it is replaced by a single plain-text string "Số giờ thực hiện".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "Số giờ thực hiện"

# Reflect the new active cell selection left by the editor after the change.
$ws.Range("D10").Select()
